# Remove the AVG / OBP / SLG stat columns (AC:AE) from the "sheet1" table.
# This mirrors the Excel UI action of selecting the three whole columns and
# choosing Delete, which shifts the trailing VIDEO column from AF back to AC.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# Select the columns that are about to be removed (matches the leftover
# selection state Excel leaves behind after a "delete entire column" edit).
$ws.Range("AC1:AE1048576").Select()

# Delete the three whole columns; remaining columns (incl. VIDEO) shift left.
$ws.Range("AC1:AE13").EntireColumn.Delete()

# Shrink the table definition down to the new extent.
$tbl.Resize($ws.Range("A1:AC13"))

# Touch the trailing header cell so the table's column-name cache picks up
# "VIDEO" (now in column AC) instead of the stale "AVG" that used to live
# there before the delete.
$ws.Range("AC1").Value = "VIDEO"
